{"js": "// Insert \"receta/\" right before \"administraci\u00f3n\" in the paragraph that\n// talks about \"Para cada administraci\u00f3n de medicamentos, ...\".\nconst body = context.document.body;\nconst results = body.search(\"administraci\u00f3n de medicamentos\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Target text \"administraci\u00f3n de medicamentos\" not found.');\n}\n\nconst found = results.items[0];\n\n// Insert \"receta/\" right before the matched \"administraci\u00f3n\" word so the\n// sentence reads \"...receta/administraci\u00f3n de medicamentos...\".\nfound.insertText(\"receta/\", \"Start\");\nawait context.sync();\n", "ps1": "# Insert \"receta/\" right before \"administraci\u00f3n\" in the paragraph that\n# talks about \"Para cada administraci\u00f3n de medicamentos, ...\".\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$found = $range.Find.Execute(\n    \"administraci\u00f3n de medicamentos\",  # FindText\n    $false,                             # MatchCase\n    $false,                             # MatchWholeWord\n    $false,                             # MatchWildcards\n    $false,                             # MatchSoundsLike\n    $false,                             # MatchAllWordForms\n    $true,                              # Forward\n    1,                                  # Wrap (wdFindContinue)\n    $false,                             # Format\n    \"receta/administraci\u00f3n de medicamentos\",  # ReplaceWith\n    2                                   # Replace (wdReplaceAll)\n)\n\nif (-not $found) {\n    throw 'Target text \"administraci\u00f3n de medicamentos\" not found.'\n}\n"}
